$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value2 = '64.443.44'
$ws.Range("E2").Value2 = '  +1.55%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value2 = '3.163.98'
$ws.Range("E3").Value2 = '  +2.18%  '
$ws.Range("E4").Value2 = '  +0.14%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value2 = '592.58'
$ws.Range("E5").Value2 = '  +1.68%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value2 = '148.19'
$ws.Range("E6").Value2 = '  +2.04%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value2 = '3.152.72'
$ws.Range("E8").Value2 = '  +2.08%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value2 = '0.534'
$ws.Range("E9").Value2 = '  +1.31%  '
$ws.Range("E10").Value2 = '  +1.71%  '
$ws.Range("E11").Value2 = '  +5.01%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value2 = '0.459'
$ws.Range("E12").Value2 = '  +0.84%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value2 = '0.0000248'
$ws.Range("E13").Value2 = '  +1.26%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value2 = '37.62'
$ws.Range("E14").Value2 = '  +1.09%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value2 = '3.676.56'
$ws.Range("E15").Value2 = '  +1.77%  '
$ws.Range("E16").Value2 = '  -0.01%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value2 = '7.29'
$ws.Range("E17").Value2 = '  +3.01%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value2 = '64.181.53'
$ws.Range("E18").Value2 = '  +1.25%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value2 = '3.151.78'
$ws.Range("E19").Value2 = '  +1.90%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value2 = '469.80'
$ws.Range("E20").Value2 = '  +1.98%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value2 = '14.55'
$ws.Range("E21").Value2 = '  +2.39%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value2 = '0.737'
$ws.Range("E22").Value2 = '  +2.15%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value2 = '7.65'
$ws.Range("E23").Value2 = '  +3.08%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value2 = '2.44'
$ws.Range("E24").Value2 = '  +14.74%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value2 = '13.21'
$ws.Range("E25").Value2 = '  +2.77%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value2 = '81.43'
$ws.Range("E26").Value2 = '  +0.28%  '
$ws.Range("B27").Value2 = 'RenderToken'
$ws.Range("C27").Value2 = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value2 = '10.10'
$ws.Range("E27").Value2 = '  +12.55%  '
$ws.Range("B28").Value2 = 'Dai'
$ws.Range("C28").Value2 = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value2 = '1.00'
$ws.Range("E28").Value2 = '  +0.12%  '
$ws.Range("E29").Value2 = '  +2.27%  '
$ws.Range("B30").Value2 = 'ImmutableX'
$ws.Range("C30").Value2 = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value2 = '2.22'
$ws.Range("E30").Value2 = '  +1.53%  '
$ws.Range("B31").Value2 = 'FirstDigitalUSD'
$ws.Range("C31").Value2 = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value2 = '1.00'
$ws.Range("E31").Value2 = '  +0.10%  '
$ws.Range("B32").Value2 = 'NEARProtocol'
$ws.Range("C32").Value2 = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value2 = '7.28'
$ws.Range("E32").Value2 = '  +5.83%  '
$ws.Range("B33").Value2 = 'Hedera'
$ws.Range("C33").Value2 = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value2 = '0.115'
$ws.Range("E33").Value2 = '  +4.73%  '
$ws.Range("B34").Value2 = 'EthereumClassic'
$ws.Range("C34").Value2 = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value2 = '27.96'
$ws.Range("E34").Value2 = '  +5.13%  '
$ws.Range("E35").Value2 = '  +2.00%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value2 = '1.07'
$ws.Range("E36").Value2 = '  +3.64%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value2 = '6.20'
$ws.Range("E37").Value2 = '  +3.78%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value2 = '2.31'
$ws.Range("E38").Value2 = '  +0.93%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value2 = '3.29'
$ws.Range("E39").Value2 = '  -3.46%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value2 = '465.98'
$ws.Range("E40").Value2 = '  +7.75%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value2 = '51.40'
$ws.Range("E41").Value2 = '  +2.31%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value2 = '9.31'
$ws.Range("E42").Value2 = '  +7.45%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value2 = '0.296'
$ws.Range("E43").Value2 = '  +9.99%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value2 = '0.0375'
$ws.Range("E44").Value2 = '  +2.59%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value2 = '2.899.80'
$ws.Range("E45").Value2 = '  +0.90%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value2 = '39.91'
$ws.Range("E46").Value2 = '  +12.59%  '
$ws.Range("E47").Value2 = '  +0.25%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value2 = '133.10'
$ws.Range("E48").Value2 = '  +7.38%  '
$ws.Range("B50").Value2 = 'ThetaToken'
$ws.Range("C50").Value2 = 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value2 = '2.24'
$ws.Range("E50").Value2 = '  +4.81%  '
$ws.Range("B51").Value2 = 'Stellar'
$ws.Range("C51").Value2 = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value2 = '0.111'
$ws.Range("E51").Value2 = '  +1.11%  '
